$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Delete the entire row 49 (Caso 6330, REPUBLICA DE LA INDIA 3106).
# This shifts rows 50-53 up to 49-52, matching the target state and
# shrinking the used range from A1:P53 to A1:P52.
$ws.Rows.Item(49).Delete()
